$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Adm"
$ws.Cells.Item(2,3).Value = "Calcr"
$ws.Cells.Item(2,4).Value = "sCs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 8.299916333333334
$ws.Cells.Item(2,8).Value = 24.899749
$ws.Cells.Item(2,9).Value = 0.2999832381170111
$ws.Cells.Item(2,10).Value = 0.2999832381170111
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 6.618147999999999
$ws.Cells.Item(2,14).Value = 19.854444
$ws.Cells.Item(2,15).Value = 1
$ws.Cells.Item(2,16).Value = 1
$ws.Cells.Item(2,17).Value = 54.93007468161733
$ws.Cells.Item(2,18).Value = 494.3706721345559
$ws.Cells.Item(2,19).Value = 0.2999832381170111
$ws.Cells.Item(2,20).Value = 0.2999832381170111
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Adm"
$ws.Cells.Item(3,3).Value = "Calcr"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 18.11265066666666
$ws.Cells.Item(3,8).Value = 54.33795199999999
$ws.Cells.Item(3,9).Value = 0.6546441409351844
$ws.Cells.Item(3,10).Value = 0.6546441409351845
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.618147999999999
$ws.Cells.Item(3,14).Value = 19.854444
$ws.Cells.Item(3,15).Value = 1
$ws.Cells.Item(3,16).Value = 1
$ws.Cells.Item(3,17).Value = 119.8722027842986
$ws.Cells.Item(3,18).Value = 1078.849825058688
$ws.Cells.Item(3,19).Value = 0.6546441409351844
$ws.Cells.Item(3,20).Value = 0.6546441409351845
$ws.Cells.Item(4,1).Value = "M2"
$ws.Cells.Item(4,2).Value = "Adm"
$ws.Cells.Item(4,3).Value = "Calcr"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.061538
$ws.Cells.Item(4,8).Value = 0.184614
$ws.Cells.Item(4,9).Value = 0.002224163204285067
$ws.Cells.Item(4,10).Value = 0.002224163204285067
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 6.618147999999999
$ws.Cells.Item(4,14).Value = 19.854444
$ws.Cells.Item(4,15).Value = 1
$ws.Cells.Item(4,16).Value = 1
$ws.Cells.Item(4,17).Value = 0.4072675916239999
$ws.Cells.Item(4,18).Value = 3.665408324616
$ws.Cells.Item(4,19).Value = 0.002224163204285067
$ws.Cells.Item(4,20).Value = 0.002224163204285067
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Adm"
$ws.Cells.Item(5,3).Value = "Calcr"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.193828666666667
$ws.Cells.Item(5,8).Value = 3.581486
$ws.Cells.Item(5,9).Value = 0.04314845774351948
$ws.Cells.Item(5,10).Value = 0.04314845774351949
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 6.618147999999999
$ws.Cells.Item(5,14).Value = 19.854444
$ws.Cells.Item(5,15).Value = 1
$ws.Cells.Item(5,16).Value = 1
$ws.Cells.Item(5,17).Value = 7.900934802642666
$ws.Cells.Item(5,18).Value = 71.10841322378398
$ws.Cells.Item(5,19).Value = 0.04314845774351948
$ws.Cells.Item(5,20).Value = 0.04314845774351949
